$d = $word.ActiveDocument

# Update the date line (first paragraph, centered, Arial 30)
$d.Content.Find.Execute("2025-03-18 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-19 Wednesday", 2) | Out-Null

# Update each cell of the single table (20 rows x 5 columns) positionally,
# since some expressions repeat (e.g. "52+41=93" appears twice) and must be
# matched by position rather than by text content.
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text = "9-0=9"
$t.Cell(1,2).Range.Text = "91+7=98"
$t.Cell(1,3).Range.Text = "55+16=71"
$t.Cell(1,4).Range.Text = "68+16=84"
$t.Cell(1,5).Range.Text = "65+17=82"

$t.Cell(2,1).Range.Text = "72-58=14"
$t.Cell(2,2).Range.Text = "52+35=87"
$t.Cell(2,3).Range.Text = "74-23=51"
$t.Cell(2,4).Range.Text = "55+6=61"
$t.Cell(2,5).Range.Text = "44+32=76"

$t.Cell(3,1).Range.Text = "95-37=58"
$t.Cell(3,2).Range.Text = "4+26=30"
$t.Cell(3,3).Range.Text = "45-17=28"
$t.Cell(3,4).Range.Text = "87-81=6"
$t.Cell(3,5).Range.Text = "49-37=12"

$t.Cell(4,1).Range.Text = "25+5=30"
$t.Cell(4,2).Range.Text = "60+25=85"
$t.Cell(4,3).Range.Text = "22+8=30"
$t.Cell(4,4).Range.Text = "72+12=84"
$t.Cell(4,5).Range.Text = "28+40=68"

$t.Cell(5,1).Range.Text = "32-1=31"
$t.Cell(5,2).Range.Text = "41-4=37"
$t.Cell(5,3).Range.Text = "12+87=99"
$t.Cell(5,4).Range.Text = "53-12=41"
$t.Cell(5,5).Range.Text = "65-59=6"

$t.Cell(6,1).Range.Text = "84-13=71"
$t.Cell(6,2).Range.Text = "2+72=74"
$t.Cell(6,3).Range.Text = "97-57=40"
$t.Cell(6,4).Range.Text = "72+4=76"
$t.Cell(6,5).Range.Text = "80-24=56"

$t.Cell(7,1).Range.Text = "19+25=44"
$t.Cell(7,2).Range.Text = "0+49=49"
$t.Cell(7,3).Range.Text = "2+11=13"
$t.Cell(7,4).Range.Text = "44+1=45"
$t.Cell(7,5).Range.Text = "73+15=88"

$t.Cell(8,1).Range.Text = "85-69=16"
$t.Cell(8,2).Range.Text = "45+45=90"
$t.Cell(8,3).Range.Text = "23+41=64"
$t.Cell(8,4).Range.Text = "15+24=39"
$t.Cell(8,5).Range.Text = "3+22=25"

$t.Cell(9,1).Range.Text = "59-8=51"
$t.Cell(9,2).Range.Text = "49+40=89"
$t.Cell(9,3).Range.Text = "70+20=90"
$t.Cell(9,4).Range.Text = "82-41=41"
$t.Cell(9,5).Range.Text = "75+8=83"

$t.Cell(10,1).Range.Text = "15+13=28"
$t.Cell(10,2).Range.Text = "84+3=87"
$t.Cell(10,3).Range.Text = "37-32=5"
$t.Cell(10,4).Range.Text = "83-3=80"
$t.Cell(10,5).Range.Text = "72-65=7"

$t.Cell(11,1).Range.Text = "77-52=25"
$t.Cell(11,2).Range.Text = "0+25=25"
$t.Cell(11,3).Range.Text = "8+61=69"
$t.Cell(11,4).Range.Text = "84+0=84"
$t.Cell(11,5).Range.Text = "58-53=5"

$t.Cell(12,1).Range.Text = "23+61=84"
$t.Cell(12,2).Range.Text = "4-4=0"
$t.Cell(12,3).Range.Text = "46-0=46"
$t.Cell(12,4).Range.Text = "76-73=3"
$t.Cell(12,5).Range.Text = "31+40=71"

$t.Cell(13,1).Range.Text = "39+9=48"
$t.Cell(13,2).Range.Text = "66-3=63"
$t.Cell(13,3).Range.Text = "67-26=41"
$t.Cell(13,4).Range.Text = "49+4=53"
$t.Cell(13,5).Range.Text = "47-46=1"

$t.Cell(14,1).Range.Text = "89-4=85"
$t.Cell(14,2).Range.Text = "37+20=57"
$t.Cell(14,3).Range.Text = "11+1=12"
$t.Cell(14,4).Range.Text = "77+5=82"
$t.Cell(14,5).Range.Text = "22+62=84"

$t.Cell(15,1).Range.Text = "44+44=88"
$t.Cell(15,2).Range.Text = "30+41=71"
$t.Cell(15,3).Range.Text = "8+59=67"
$t.Cell(15,4).Range.Text = "3+13=16"
$t.Cell(15,5).Range.Text = "47+14=61"

$t.Cell(16,1).Range.Text = "45-5=40"
$t.Cell(16,2).Range.Text = "26+16=42"
$t.Cell(16,3).Range.Text = "7+2=9"
$t.Cell(16,4).Range.Text = "68-52=16"
$t.Cell(16,5).Range.Text = "81-32=49"

$t.Cell(17,1).Range.Text = "84-47=37"
$t.Cell(17,2).Range.Text = "21+39=60"
$t.Cell(17,3).Range.Text = "62+0=62"
$t.Cell(17,4).Range.Text = "29+67=96"
$t.Cell(17,5).Range.Text = "54-2=52"

$t.Cell(18,1).Range.Text = "30+13=43"
$t.Cell(18,2).Range.Text = "67+4=71"
$t.Cell(18,3).Range.Text = "42+10=52"
$t.Cell(18,4).Range.Text = "39+23=62"
$t.Cell(18,5).Range.Text = "31+5=36"

$t.Cell(19,1).Range.Text = "78-51=27"
$t.Cell(19,2).Range.Text = "43-10=33"
$t.Cell(19,3).Range.Text = "14+84=98"
$t.Cell(19,4).Range.Text = "38+49=87"
$t.Cell(19,5).Range.Text = "26+46=72"

$t.Cell(20,1).Range.Text = "83-61=22"
$t.Cell(20,2).Range.Text = "68-21=47"
$t.Cell(20,3).Range.Text = "60+30=90"
$t.Cell(20,4).Range.Text = "84+5=89"
$t.Cell(20,5).Range.Text = "57-7=50"
